# Updated symbol list on Tue Dec 27 16:34:31 UTC 2022 with GitHub Actions
#
# This script reproduces the per-row price/volume refresh (and the
# cyclic re-ranking of rows 18-24) recorded in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text. Cells in this sheet store everything
# (including numbers) as text, so purely-numeric-looking strings need a
# leading apostrophe to stop Excel from auto-converting them to the
# Number type.
function Set-TextValue($cell, [string]$value) {
    if ($value -match '^-?\d+(\.\d+)?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}

# ---------------------------------------------------------------------
# Simple price (column D) refreshes that do not touch the coin order
# ---------------------------------------------------------------------
$priceUpdates = @{
    3  = "23.70"
    4  = "5.369"
    5  = "0.05879"
    6  = "3.374"
    7  = "6.480"
    8  = "0.8112"
    9  = "0.9227"
    10 = "0.1416"
    11 = "0.07389"
    12 = "0.03088"
    13 = "0.03056"
    14 = "0.09354"
    15 = "3.873"
    16 = "0.001559"
    17 = "0.04691"
    25 = "0.3231"
    26 = "0.1330"
    40 = "0.03876"
    41 = "0.006370"
    43 = "0.003200"
    44 = "0.008574"
    45 = "0.00005250"
    47 = "0.6811"
    48 = "0.001676"
}

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $ws.Cells.Item($row, 4) $priceUpdates[$row]
}

# ---------------------------------------------------------------------
# Rows 18-24: the coin rankings shifted by one place (row 18's coin
# moved to 24, everything else moved up one row) and the price/volume
# columns were refreshed with new data.
# ---------------------------------------------------------------------
$rowData = @{
    18 = @("TigerCash",   "https://coinranking.com/coin/6hIn06L2+tigercash-tch",              "0.005955",   "17TigerCashTCH")
    19 = @("BitKan",      "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan",             "0.001241",   "18BitKanKAN")
    20 = @("HotbitToken", "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb",         "0.004720",   "19HotbitTokenHTB")
    21 = @("NitroEx",     "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx",              "0.00008801", "20NitroExNTX")
    22 = @("LEO",         "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo",                 "3.596",      "21LEOLEO")
    23 = @("BTSEToken",   "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse",          "2.158",      "22BTSETokenBTSE")
    24 = @("One",         "https://coinranking.com/coin/6Lga5NiXX3rT+one-one",                 "0.01094",    "23OneONEBestin24h")
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    Set-TextValue $ws.Cells.Item($row, 2) $vals[0]
    Set-TextValue $ws.Cells.Item($row, 3) $vals[1]
    Set-TextValue $ws.Cells.Item($row, 4) $vals[2]
    Set-TextValue $ws.Cells.Item($row, 5) $vals[3]
}
